$d = $word.ActiveDocument

# Original paragraph text: "Version 2." split across runs as
# "Versi" | "on" | [spellEnd] | " 2" | [bookmark] | "."
# Target: "Version" | [spellEnd] | " 1." | [bookmark]  (no trailing "." run)
#
# Work from the tail of the range toward the front so earlier offsets
# stay valid as later edits shrink/grow the text.

# 1) Remove the standalone trailing "." run (position 9-10).
$d.Range(9, 10).Delete()

# 2) Turn " 2" (position 7-9) into " 1." (keeps it inside the same run).
$d.Range(7, 9).Text = " 1."

# 3) Merge "Versi" + "on" into a single "Version" run: delete the "on"
#    text (position 5-7) and re-insert it right after "Versi" so it
#    collapses back into the first run instead of staying split.
$d.Range(5, 7).Delete()
$d.Range(5, 5).InsertAfter("on")
